# Update Name of Algo
# Apply corrected imputed values produced by the KNN algorithm run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value  = -13.376
$ws.Range("A9").Value  = -21.534
$ws.Range("C9").Value  = -10.52
$ws.Range("C11").Value = -12.601
$ws.Range("A18").Value = -21.898
$ws.Range("A20").Value = -20.623
$ws.Range("C23").Value = -13.046
$ws.Range("C24").Value = -12.38
$ws.Range("C26").Value = -12.696
$ws.Range("A27").Value = -21.868
$ws.Range("C34").Value = -12.049
$ws.Range("A35").Value = -20.186
$ws.Range("C35").Value = -12.12
$ws.Range("C48").Value = -11.732
$ws.Range("C49").Value = -13.285
$ws.Range("C52").Value = -11.665
$ws.Range("C66").Value = -11.574
$ws.Range("C67").Value = -10.875
$ws.Range("A69").Value = -21.524
$ws.Range("A76").Value = -20.241
$ws.Range("A78").Value = -20.086
$ws.Range("C78").Value = -12.924
$ws.Range("C80").Value = -12.664
$ws.Range("A82").Value = -22.004
$ws.Range("A83").Value = -21.84
$ws.Range("A93").Value = -21.475
$ws.Range("C99").Value = -11.95
$ws.Range("C104").Value = -12.822
